# Apply the "api onboarding and check module file" edit:
#  1. Add a new "data_freq" column (K) to the "plant" sheet with a value per row.
#  2. Move the selection on "plant" to L3 (and it is no longer the active tab).
#  3. Mark the "plant" sheet's print orientation as portrait (adds <pageSetup>).
#  4. Insert a new "demo" worksheet between "plant" and "asset", seeded with the
#     header row plus the first two data rows from "plant" (renamed to "Demo Project").

$wb = $excel.ActiveWorkbook
$plant = $wb.Worksheets.Item("plant")

# --- 1. data_freq header + per-row values on "plant" ---------------------
$j1 = $plant.Range("J1")
$k1 = $plant.Range("K1")
$j1.Copy($k1)
$k1.Value = "data_freq"

$freq = @{
    2 = 1; 3 = 5; 4 = 5; 5 = 5; 6 = 5; 7 = 5; 8 = 5; 9 = 5; 10 = 5;
    11 = 5; 12 = 5; 13 = 5; 14 = 5; 15 = 5; 16 = 5; 17 = 5; 18 = 5; 19 = 5; 20 = 5;
    21 = 5; 22 = 5; 23 = 5; 24 = 5; 25 = 5; 26 = 5; 27 = 5; 28 = 5; 29 = 5; 30 = 5;
    31 = 5; 32 = 5; 33 = 5; 34 = 5; 35 = 5; 36 = 5; 37 = 5; 38 = 5; 39 = 5; 40 = 5;
    41 = 5; 42 = 5; 43 = 5; 44 = 5; 45 = 5; 46 = 5; 47 = 5; 48 = 5; 49 = 5; 50 = 5;
    51 = 5; 52 = 5; 53 = 5; 54 = 5; 55 = 5; 56 = 5; 57 = 5; 58 = 5; 59 = 5; 60 = 5;
    61 = 5; 62 = 5; 63 = 5; 64 = 5; 65 = 5; 66 = 5; 67 = 5; 68 = 5; 69 = 5; 70 = 5;
    71 = 5; 72 = 5; 73 = 5; 74 = 5; 75 = 5; 76 = 5; 77 = 5; 78 = 5; 79 = 5; 80 = 5;
    81 = 5; 82 = 5; 83 = 5; 84 = 5; 85 = 1; 86 = 5; 87 = 5; 88 = 5; 89 = 5; 90 = 5;
    91 = 10
}
foreach ($row in $freq.Keys) {
    $plant.Cells.Item($row, 11).Value = $freq[$row]
}

# --- 2. selection moves to L3 on "plant" ----------------------------------
$plant.Range("L3").Select()

# --- 3. print orientation recorded on "plant" -----------------------------
$plant.PageSetup.Orientation = 1

# --- 4. new "demo" sheet between "plant" and "asset" ----------------------
$demo = $wb.Worksheets.Add($null, $plant)
$demo.Name = "demo"

# Header row: clone formatting + shared strings straight from "plant" row 1.
$plant.Range("A1:J1").Copy($demo.Range("A1:J1"))
$plant.Range("K1").Copy($demo.Range("K1"))

# Data rows: clone formatting/values from "plant" rows 2-3, then rename the project.
$plant.Range("A2:J2").Copy($demo.Range("A2:J2"))
$plant.Range("A3:J3").Copy($demo.Range("A3:J3"))
$demo.Range("C2").Value = "Demo Project"
$demo.Range("C3").Value = "Demo Project"

$demo.Range("C8").Select()

Write-Output "edit applied"
